$d = $word.ActiveDocument

# Update the date/title line (first paragraph)
$d.Paragraphs.Item(1).Range.Text = '2024-07-04 Thursday'

# Update each table cell by directly assigning the new text to its Range.
# (Find/Replace scoped to a cell Range incorrectly matched duplicate text
# elsewhere in the document, e.g. the two "32-0=" cells, so we set .Text directly
# instead, which only touches the targeted cell and preserves its run formatting.)
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = '76-8='
$t.Cell(1, 2).Range.Text = '76+9='
$t.Cell(1, 3).Range.Text = '30+63='
$t.Cell(1, 4).Range.Text = '93-65='
$t.Cell(1, 5).Range.Text = '33+8='
$t.Cell(2, 1).Range.Text = '9-4='
$t.Cell(2, 2).Range.Text = '31+42='
$t.Cell(2, 3).Range.Text = '52+12='
$t.Cell(2, 4).Range.Text = '13+72='
$t.Cell(2, 5).Range.Text = '47+52='
$t.Cell(3, 1).Range.Text = '51-49='
$t.Cell(3, 2).Range.Text = '97-86='
$t.Cell(3, 3).Range.Text = '1+7='
$t.Cell(3, 4).Range.Text = '80-18='
$t.Cell(3, 5).Range.Text = '8+39='
$t.Cell(4, 1).Range.Text = '59-31='
$t.Cell(4, 2).Range.Text = '39+39='
$t.Cell(4, 3).Range.Text = '68-13='
$t.Cell(4, 4).Range.Text = '16+19='
$t.Cell(4, 5).Range.Text = '79+14='
$t.Cell(5, 1).Range.Text = '14+46='
$t.Cell(5, 2).Range.Text = '75-4='
$t.Cell(5, 3).Range.Text = '56-46='
$t.Cell(5, 4).Range.Text = '94-6='
$t.Cell(5, 5).Range.Text = '30+46='
$t.Cell(6, 1).Range.Text = '20+8='
$t.Cell(6, 2).Range.Text = '14+5='
$t.Cell(6, 3).Range.Text = '26+63='
$t.Cell(6, 4).Range.Text = '59-49='
$t.Cell(6, 5).Range.Text = '28+70='
$t.Cell(7, 1).Range.Text = '19+37='
$t.Cell(7, 2).Range.Text = '36-10='
$t.Cell(7, 3).Range.Text = '25+10='
$t.Cell(7, 4).Range.Text = '19+64='
$t.Cell(7, 5).Range.Text = '12-10='
$t.Cell(8, 1).Range.Text = '89-79='
$t.Cell(8, 2).Range.Text = '37+19='
$t.Cell(8, 3).Range.Text = '92-15='
$t.Cell(8, 4).Range.Text = '33+34='
$t.Cell(8, 5).Range.Text = '89-61='
$t.Cell(9, 1).Range.Text = '95-78='
$t.Cell(9, 2).Range.Text = '76-16='
$t.Cell(9, 3).Range.Text = '51+2='
$t.Cell(9, 4).Range.Text = '64-18='
$t.Cell(9, 5).Range.Text = '43+9='
$t.Cell(10, 1).Range.Text = '44-44='
$t.Cell(10, 2).Range.Text = '10+65='
$t.Cell(10, 3).Range.Text = '38-15='
$t.Cell(10, 4).Range.Text = '9+83='
$t.Cell(10, 5).Range.Text = '31-20='
$t.Cell(11, 1).Range.Text = '50+24='
$t.Cell(11, 2).Range.Text = '64+28='
$t.Cell(11, 3).Range.Text = '45+35='
$t.Cell(11, 4).Range.Text = '72+13='
$t.Cell(11, 5).Range.Text = '73-13='
$t.Cell(12, 1).Range.Text = '87-85='
$t.Cell(12, 2).Range.Text = '45+45='
$t.Cell(12, 3).Range.Text = '89-68='
$t.Cell(12, 4).Range.Text = '30-27='
$t.Cell(12, 5).Range.Text = '85-5='
$t.Cell(13, 1).Range.Text = '61-4='
$t.Cell(13, 2).Range.Text = '2+21='
$t.Cell(13, 3).Range.Text = '33+43='
$t.Cell(13, 4).Range.Text = '83-36='
$t.Cell(13, 5).Range.Text = '48-27='
$t.Cell(14, 1).Range.Text = '62+35='
$t.Cell(14, 2).Range.Text = '89-47='
$t.Cell(14, 3).Range.Text = '58+12='
$t.Cell(14, 4).Range.Text = '24-10='
$t.Cell(14, 5).Range.Text = '76-5='
$t.Cell(15, 1).Range.Text = '10+62='
$t.Cell(15, 2).Range.Text = '30+67='
$t.Cell(15, 3).Range.Text = '33+40='
$t.Cell(15, 4).Range.Text = '27+33='
$t.Cell(15, 5).Range.Text = '24+52='
$t.Cell(16, 1).Range.Text = '72-14='
$t.Cell(16, 2).Range.Text = '87-79='
$t.Cell(16, 3).Range.Text = '8+27='
$t.Cell(16, 4).Range.Text = '72-22='
$t.Cell(16, 5).Range.Text = '41-34='
$t.Cell(17, 1).Range.Text = '99-60='
$t.Cell(17, 2).Range.Text = '30-3='
$t.Cell(17, 3).Range.Text = '77-16='
$t.Cell(17, 4).Range.Text = '60-31='
$t.Cell(17, 5).Range.Text = '4+50='
$t.Cell(18, 1).Range.Text = '74-63='
$t.Cell(18, 2).Range.Text = '38+60='
$t.Cell(18, 3).Range.Text = '23+41='
$t.Cell(18, 4).Range.Text = '37+27='
$t.Cell(18, 5).Range.Text = '91-62='
$t.Cell(19, 1).Range.Text = '19-11='
$t.Cell(19, 2).Range.Text = '77+11='
$t.Cell(19, 3).Range.Text = '84-63='
$t.Cell(19, 4).Range.Text = '41+14='
$t.Cell(19, 5).Range.Text = '87-32='
$t.Cell(20, 1).Range.Text = '35+27='
$t.Cell(20, 2).Range.Text = '71-47='
$t.Cell(20, 3).Range.Text = '86-0='
$t.Cell(20, 4).Range.Text = '60+30='
$t.Cell(20, 5).Range.Text = '38+37='
